$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores numeric-looking values as plain text (e.g. "72.346.07",
# "0.171", "1.00") because some of them use "." as a thousands separator, which is not
# a valid Excel number. Apply a Text number format across the whole column first so
# that assigning these strings below does not get auto-converted into real numbers
# (which would silently drop meaningful trailing/insignificant zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "72.209.02"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.630.90"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "582.21"
$ws.Range("E5").Value = "  -3.64%  "
$ws.Range("D6").Value = "173.64"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").Value = "2.630.02"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "0.354"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "4.92"
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").Value = "3.119.08"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "72.110.57"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000184"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").Value = "25.76"
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").Value = "2.622.97"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").Value = "12.09"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "7.88"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "372.95"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").Value = "4.10"
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").Value = "2.05"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "70.77"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").Value = "4.20"
$ws.Range("E26").Value = "  -4.00%  "
$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  -4.11%  "
$ws.Range("D28").Value = "2.770.73"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "0.0₃0945"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "7.92"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("D32").Value = "493.64"
$ws.Range("E32").Value = "  -5.48%  "
$ws.Range("D33").Value = "1.27"
$ws.Range("E33").Value = "  -3.12%  "
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "163.02"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "19.17"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("D39").Value = "18.84"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "1.71"
$ws.Range("E42").Value = "  -7.62%  "
$ws.Range("D43").Value = "2.54"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").Value = "4.87"
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("D45").Value = "0.324"
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("D46").Value = "39.02"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "151.85"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "3.63"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("D49").Value = "0.542"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").Value = "1.66"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").Value = "0.598"
$ws.Range("E51").Value = "  -0.75%  "
